$d = $word.ActiveDocument
$findText = "significance level."
$replaceText = "significance level.^p" + 'Answer)^p^pX <- data.frame(UCBAdmissions)^pAdmitted <- X$Admit^pFreq <- X$Freq^p^p# From the question, the value of p0 is 0.4. ^p# Sample proportion = (# of admitted)/Total^pp0 = 0.4^pAdm <- grep("Admitted", Admitted)^pAdm^pAdm1 = NULL^pfor(i in Adm)^p{^p X <- Freq[i]^p Adm1 = rbind(Adm1, X)^p}^p^pno.adm <- sum(Adm1)   # number of students admitted ^pTotal <- sum(Freq) # total students^pp_hat = no.adm / Total  # Sample proportion^p^p# Test statistic^pz = (p_hat - p0) / round(sqrt((p0*(1-p0))/Total),8)^pz^pprob = pnorm(z, 0, 1)^pprob^p^p# The probability value obtained is 0.04638927, which is greater than 0.01. Therefore^p# we can conclude that the acceptance range is consistent with the officer''s claim.'
$found = $d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
Write-Output "Found: $found"
Write-Output "Paragraph count: $($d.Paragraphs.Count)"
